# Update on 12 Nov 2017
# Rename the column headers in row 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Vendor Name"
$ws.Range("C1").Value = "Location From"
$ws.Range("D1").Value = "Location To"
$ws.Range("E1").Value = "KG"
$ws.Range("F1").Value = "Trips"

# Column F held shipment dates; replace them with trip counts (numbers).
$ws.Range("F2").Value = 5
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 139
$ws.Range("F5").Value = 9
$ws.Range("F6").Value = 39
$ws.Range("F7").Value = 4

# Columns A:B had custom widths; reset the whole sheet to a uniform width.
$ws.Range($ws.Columns.Item(1), $ws.Columns.Item(1025)).ColumnWidth = 8.57085020242915

# Tab ratio (split between sheet tabs and horizontal scrollbar) nudged slightly.
$excel.ActiveWindow.TabRatio = 0.994

# Selection moved from F9 to F8.
$ws.Range("F8").Select()
